$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 8, pushing the existing
# "2024-05-18" row (old row 8) and the dated-serial row (old row 9)
# down to rows 9 and 11 respectively, leaving rows 8 and 10 free for
# the two newly-reported days (2024-05-17 and 2024-05-19).
$ws.Rows("8:9").Insert()

# Helper: write a value into a cell while forcing text storage so that
# strings like "2024-05-17" or "3,147" are NOT auto-converted into a
# date serial / grouped number by Excel's input parser.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- New row 8: 2024-05-17 -------------------------------------------------
Set-TextCell 8 1 "2024-05-17"
Set-TextCell 8 2 "3,147"
Set-TextCell 8 3 "1,988"
Set-TextCell 8 4 "5,224"
Set-TextCell 8 5 "5,831"
Set-TextCell 8 6 "4,303"
Set-TextCell 8 7 "3,185"

# --- Row 9 (was row 8): 2024-05-18, values unchanged ------------------------
Set-TextCell 9 1 "2024-05-18"
Set-TextCell 9 2 "3,152"
Set-TextCell 9 3 "2,005"
Set-TextCell 9 4 "5,222"
Set-TextCell 9 5 "5,831"
Set-TextCell 9 6 "4,301"
Set-TextCell 9 7 "3,175"

# --- New row 10: 2024-05-19 -------------------------------------------------
Set-TextCell 10 1 "2024-05-19"
Set-TextCell 10 2 "3,153"
Set-TextCell 10 3 "2,007"
Set-TextCell 10 4 "5,222"
Set-TextCell 10 5 "5,837"
Set-TextCell 10 6 "4,301"
Set-TextCell 10 7 "3,175"

# --- Row 11 (was row 9): keeps the numeric date-serial A cell, but with
#     an updated serial value (45432) and updated B/ C values -------------
$ws.Cells.Item(11, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(11, 1).Value = 45432
Set-TextCell 11 2 "3,152"
Set-TextCell 11 3 "2,007"
Set-TextCell 11 4 "5,222"
Set-TextCell 11 5 "5,837"
Set-TextCell 11 6 "4,300"
Set-TextCell 11 7 "3,174"
